$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Temp. & Hum. Sensor (SHT30-DIS-B)
$ws.Range("B4").Value = "Temp. & Hum. Sensor"
$ws.Range("C4").Value = "`nSHT30-DIS-B2.5kS; I2C Interface; 2.15V - 5.5V Supply; Accuracy of +/-2%RH and +/-0.2°C"
$ws.Range("D4").Value = "8-Pin DFN"
$ws.Range("F4").Value = "403-SHT30-DIS-B "

# Row 5 - Protective Cover (SHT30-DIS-F)
$ws.Range("B5").Value = "Protective Cover"
$ws.Range("C5").Value = "`nSHT30-DIS-F2.5kS; Productive Cover for Temp. & Hum. Sensor"
$ws.Range("D5").Value = "-"
$ws.Range("F5").Value = "403-SHT30-DIS-F "

# Row 6 - Gyroskop (I3G4250DTR)
$ws.Range("B6").Value = "Gyroskop"
$ws.Range("C6").Value = "I3G4250DTR; X, Y, Z; 2.4V - 3.6V Supply; I2C/SPI Interface; programmable Range"
$ws.Range("D6").Value = "LGA-16"
$ws.Range("F6").Value = "511-I3G4250DTR "
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 7.18

# Row 7 - Beschleunigungssensor (BMA423)
$ws.Range("B7").Value = "Beschleunigungssensor"
$ws.Range("C7").Value = "BMA423; X, Y, Z; 1.2V - 3.6V Supply; I2C/SPI Interface; programmable Range"
$ws.Range("D7").Value = "LGA-12"
$ws.Range("F7").Value = "262-BMA423 "
$ws.Range("H7").Value = 2.01

# Row 8 - Ultraschall-Distanzmesser (HC-SR04)
$ws.Range("B8").Value = "Ultraschall-Distanzmesser"
$ws.Range("C8").Value = "`nHC-SR04; PWM Out; 5V Supply"
$ws.Range("D8").Value = "-"
$ws.Range("F8").Value = "375-HC-SR04 "
$ws.Range("H8").Value = 5.44

# Rows 9-17: clear the old BOM entries (B..H); formulas in column I remain
$ws.Range("B9:H17").ClearContents()

# Row 22: restore the shared total formula
$ws.Range("I22").Formula = "=G22*H22"

# Row heights (content rewrap after the BOM rewrite)
$ws.Rows("4").RowHeight = 82.5
$ws.Rows("5").RowHeight = 66
$ws.Rows("6").RowHeight = 66
$ws.Rows("7").RowHeight = 49.5
$ws.Rows("8").RowHeight = 49.5
$ws.Rows("11").RowHeight = 16.5
$ws.Rows("12").RowHeight = 16.5
$ws.Rows("13").RowHeight = 16.5
$ws.Rows("17").RowHeight = 16.5

# Selection moves to F11
$ws.Range("F11").Select() | Out-Null
